# Weekly update: insert this week's 3 new "Chirimoya" price rows
# (Macroferia Regional de Talca, date 44839) right before the existing
# row 74, pushing the rest of the table down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 74-76 (rows 74..101 shift down to 77..104).
$ws.Range("A74:A76").EntireRow.Insert()

# New row 74: Especial
$ws.Cells.Item(74, 1).Value = 5
$ws.Cells.Item(74, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value = "Maule"
$ws.Cells.Item(74, 4).Value = 44839
$ws.Cells.Item(74, 5).Value = 7
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100107
$ws.Cells.Item(74, 8).Value = "Otros"
$ws.Cells.Item(74, 9).Value = 100107002
$ws.Cells.Item(74, 10).Value = "Chirimoya"
$ws.Cells.Item(74, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(74, 12).Value = "Especial"
$ws.Cells.Item(74, 13).Value = 40
$ws.Cells.Item(74, 14).Value = 25000
$ws.Cells.Item(74, 15).Value = 25000
$ws.Cells.Item(74, 16).Value = 25000
$ws.Cells.Item(74, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(74, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(74, 19).Value = 2500
$ws.Cells.Item(74, 20).Value = 10

# New row 75: Primera
$ws.Cells.Item(75, 1).Value = 5
$ws.Cells.Item(75, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(75, 3).Value = "Maule"
$ws.Cells.Item(75, 4).Value = 44839
$ws.Cells.Item(75, 5).Value = 7
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100107
$ws.Cells.Item(75, 8).Value = "Otros"
$ws.Cells.Item(75, 9).Value = 100107002
$ws.Cells.Item(75, 10).Value = "Chirimoya"
$ws.Cells.Item(75, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 50
$ws.Cells.Item(75, 14).Value = 22000
$ws.Cells.Item(75, 15).Value = 22000
$ws.Cells.Item(75, 16).Value = 22000
$ws.Cells.Item(75, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(75, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(75, 19).Value = 2200
$ws.Cells.Item(75, 20).Value = 10

# New row 76: Segunda
$ws.Cells.Item(76, 1).Value = 5
$ws.Cells.Item(76, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(76, 3).Value = "Maule"
$ws.Cells.Item(76, 4).Value = 44839
$ws.Cells.Item(76, 5).Value = 7
$ws.Cells.Item(76, 6).Value = "Fruta"
$ws.Cells.Item(76, 7).Value = 100107
$ws.Cells.Item(76, 8).Value = "Otros"
$ws.Cells.Item(76, 9).Value = 100107002
$ws.Cells.Item(76, 10).Value = "Chirimoya"
$ws.Cells.Item(76, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(76, 12).Value = "Segunda"
$ws.Cells.Item(76, 13).Value = 30
$ws.Cells.Item(76, 14).Value = 20000
$ws.Cells.Item(76, 15).Value = 20000
$ws.Cells.Item(76, 16).Value = 20000
$ws.Cells.Item(76, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(76, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(76, 19).Value = 2000
$ws.Cells.Item(76, 20).Value = 10
